$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was collected for row 86 ("Pepino ensalada" at
# Vega Monumental Concepción). Insert a fresh row at position 86, pushing
# the existing data (rows 86-149) down to rows 87-150, then populate the
# new row with the latest observation.
$ws.Rows.Item(86).Insert()

$ws.Cells.Item(86, 1).Value = 11
$ws.Cells.Item(86, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(86, 3).Value = "Bíobío"
$ws.Cells.Item(86, 4).Value = 44790
$ws.Cells.Item(86, 5).Value = 8
$ws.Cells.Item(86, 6).Value = 100112043
$ws.Cells.Item(86, 7).Value = "Pepino ensalada"
$ws.Cells.Item(86, 8).Value = "Sin especificar"
$ws.Cells.Item(86, 9).Value = "Primera"
$ws.Cells.Item(86, 10).Value = 270
$ws.Cells.Item(86, 11).Value = 20000
$ws.Cells.Item(86, 12).Value = 21000
$ws.Cells.Item(86, 13).Value = 20444
$ws.Cells.Item(86, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(86, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(86, 16).Value = 341
$ws.Cells.Item(86, 17).Value = 60
$ws.Cells.Item(86, 18).Value = "Hortaliza"
